# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (columns H-N: currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) on a handful of leve rows across the job sheets,
# mirroring the scheduled price-refresh job's committed diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 149.25
$ws.Range("J29").Value = 149.25
$ws.Range("L29").Value = 447.75
$ws.Range("N29").Value = -1009.75
$ws.Range("H38").Value = 2437.5881
$ws.Range("I38").Value = 620
$ws.Range("J38").Value = 4053.2222
$ws.Range("K38").Value = 1860
$ws.Range("L38").Value = 12159.6666
$ws.Range("M38").Value = -1488
$ws.Range("N38").Value = -12903.6666
$ws.Range("H58").Value = 51
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null
$ws.Range("H74").Value = 4067.1667
$ws.Range("I74").Value = 4067.1667
$ws.Range("K74").Value = 4067.1667
$ws.Range("M74").Value = -3131.1667
$ws.Range("H77").Value = 4067.1667
$ws.Range("I77").Value = 4067.1667
$ws.Range("K77").Value = 20335.8335
$ws.Range("M77").Value = -15655.8335
$ws.Range("H138").Value = 3589.5667
$ws.Range("J138").Value = 3816.9167
$ws.Range("L138").Value = 11450.7501
$ws.Range("N138").Value = -21730.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4526.9473
$ws.Range("I32").Value = 2276.9312
$ws.Range("K32").Value = 2276.9312
$ws.Range("M32").Value = -1989.9312
$ws.Range("H74").Value = 2620.8975
$ws.Range("J74").Value = 2677.647
$ws.Range("L74").Value = 2677.647
$ws.Range("N74").Value = -4425.647
$ws.Range("H77").Value = 2620.8975
$ws.Range("J77").Value = 2677.647
$ws.Range("L77").Value = 13388.235
$ws.Range("N77").Value = -22124.235
$ws.Range("H132").Value = 3099.25
$ws.Range("I132").Value = 3095.3953
$ws.Range("J132").Value = 3265
$ws.Range("K132").Value = 9286.1859
$ws.Range("L132").Value = 9795
$ws.Range("M132").Value = -6756.1859
$ws.Range("N132").Value = -14855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 251859.17
$ws.Range("I86").Value = 334992.4
$ws.Range("J86").Value = 2459.5
$ws.Range("K86").Value = 334992.4
$ws.Range("L86").Value = 2459.5
$ws.Range("M86").Value = -333869.4
$ws.Range("N86").Value = -4705.5
$ws.Range("H89").Value = 251859.17
$ws.Range("I89").Value = 334992.4
$ws.Range("J89").Value = 2459.5
$ws.Range("K89").Value = 1674962
$ws.Range("L89").Value = 12297.5
$ws.Range("M89").Value = -1669346
$ws.Range("N89").Value = -23529.5
$ws.Range("H94").Value = 976
$ws.Range("I94").Value = 801.75
$ws.Range("J94").Value = 1324.5
$ws.Range("K94").Value = 801.75
$ws.Range("L94").Value = 1324.5
$ws.Range("M94").Value = -350.75
$ws.Range("N94").Value = -2226.5
$ws.Range("H134").Value = 12971.48
$ws.Range("I134").Value = 5152.6875
$ws.Range("J134").Value = 26871.555
$ws.Range("K134").Value = 15458.0625
$ws.Range("L134").Value = 80614.66500000001
$ws.Range("M134").Value = -12923.0625
$ws.Range("N134").Value = -85684.66500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2587.4546
$ws.Range("I31").Value = 1157.0555
$ws.Range("J31").Value = 3283.3242
$ws.Range("K31").Value = 1157.0555
$ws.Range("L31").Value = 3283.3242
$ws.Range("M31").Value = -862.0554999999999
$ws.Range("N31").Value = -3873.3242
$ws.Range("H34").Value = 2587.4546
$ws.Range("I34").Value = 1157.0555
$ws.Range("J34").Value = 3283.3242
$ws.Range("K34").Value = 1157.0555
$ws.Range("L34").Value = 3283.3242
$ws.Range("M34").Value = -955.0554999999999
$ws.Range("N34").Value = -3687.3242
$ws.Range("H58").Value = 5156.7144
$ws.Range("I58").Value = 3504.647
$ws.Range("J58").Value = 12178
$ws.Range("K58").Value = 3504.647
$ws.Range("L58").Value = 12178
$ws.Range("M58").Value = -3301.647
$ws.Range("N58").Value = -12584
$ws.Range("H132").Value = 22474.814
$ws.Range("I132").Value = 13761.456
$ws.Range("J132").Value = 43570.316
$ws.Range("K132").Value = 41284.368
$ws.Range("L132").Value = 130710.948
$ws.Range("M132").Value = -38754.368
$ws.Range("N132").Value = -135770.948
$ws.Range("H136").Value = 5156.7144
$ws.Range("I136").Value = 3504.647
$ws.Range("J136").Value = 12178
$ws.Range("K136").Value = 10513.941
$ws.Range("L136").Value = 36534
$ws.Range("M136").Value = -7963.940999999999
$ws.Range("N136").Value = -41634

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1075.3182
$ws.Range("I5").Value = 752.6
$ws.Range("J5").Value = 1766.8572
$ws.Range("K5").Value = 2257.8
$ws.Range("L5").Value = 5300.571599999999
$ws.Range("M5").Value = -2145.8
$ws.Range("N5").Value = -5524.571599999999
$ws.Range("H46").Value = 200.5
$ws.Range("I46").Value = 203
$ws.Range("J46").Value = 198
$ws.Range("K46").Value = 609
$ws.Range("L46").Value = 594
$ws.Range("M46").Value = -518
$ws.Range("N46").Value = -776
$ws.Range("H113").Value = 1040.8
$ws.Range("I113").Value = 1554.25
$ws.Range("K113").Value = 4662.75
$ws.Range("M113").Value = -2492.75
$ws.Range("H129").Value = 2646.2307
$ws.Range("I129").Value = 1254.75
$ws.Range("J129").Value = 4872.6
$ws.Range("K129").Value = 3764.25
$ws.Range("L129").Value = 14617.8
$ws.Range("M129").Value = 1235.75
$ws.Range("N129").Value = -24617.8
$ws.Range("H131").Value = 16365.464
$ws.Range("J131").Value = 2117.0356
$ws.Range("L131").Value = 6351.1068
$ws.Range("N131").Value = -16431.1068
$ws.Range("H132").Value = 7170024
$ws.Range("J132").Value = 16727057
$ws.Range("L132").Value = 150543513
$ws.Range("N132").Value = -150548573
$ws.Range("H134").Value = 1675.3846
$ws.Range("I134").Value = 1675.3846
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5026.1538
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 43.84619999999995
$ws.Range("N134").Value = $null
$ws.Range("H135").Value = 1075.3182
$ws.Range("I135").Value = 752.6
$ws.Range("J135").Value = 1766.8572
$ws.Range("K135").Value = 6773.400000000001
$ws.Range("L135").Value = 15901.7148
$ws.Range("M135").Value = -4238.400000000001
$ws.Range("N135").Value = -20971.7148
$ws.Range("H139").Value = 1976.3334
$ws.Range("I139").Value = 1976.3334
$ws.Range("K139").Value = 5929.0002
$ws.Range("M139").Value = -789.0002000000004
$ws.Range("H140").Value = 2731.375
$ws.Range("I140").Value = 2335.8572
$ws.Range("K140").Value = 7007.571599999999
$ws.Range("M140").Value = -1827.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 57499.75
$ws.Range("J7").Value = 57499.75
$ws.Range("L7").Value = 57499.75
$ws.Range("N7").Value = -57723.75
$ws.Range("H8").Value = 57499.75
$ws.Range("J8").Value = 57499.75
$ws.Range("L8").Value = 57499.75
$ws.Range("N8").Value = -57777.75
$ws.Range("H97").Value = 869.7105
$ws.Range("I97").Value = 708.9231
$ws.Range("K97").Value = 708.9231
$ws.Range("M97").Value = -212.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 191
$ws.Range("I55").Value = 97.09999999999999
$ws.Range("J55").Value = 378.8
$ws.Range("K55").Value = 97.09999999999999
$ws.Range("L55").Value = 378.8
$ws.Range("M55").Value = 75.90000000000001
$ws.Range("N55").Value = -724.8
$ws.Range("H100").Value = 4643.5
$ws.Range("I100").Value = 3491.3333
$ws.Range("K100").Value = 3491.3333
$ws.Range("M100").Value = -2950.3333
$ws.Range("H132").Value = 3173.932
$ws.Range("I132").Value = 2771.182
$ws.Range("J132").Value = 4382.1816
$ws.Range("K132").Value = 8313.545999999998
$ws.Range("L132").Value = 13146.5448
$ws.Range("M132").Value = -5783.545999999998
$ws.Range("N132").Value = -18206.5448
$ws.Range("H136").Value = 2037.5
$ws.Range("I136").Value = 1864.7059
$ws.Range("K136").Value = 5594.1177
$ws.Range("M136").Value = -3044.1177

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null
$ws.Range("H46").Value = 69928.5
$ws.Range("J46").Value = 69928.5
$ws.Range("L46").Value = 69928.5
$ws.Range("N46").Value = -70390.5
$ws.Range("H132").Value = 20335.883
$ws.Range("I132").Value = 12831.8
$ws.Range("J132").Value = 35344.05
$ws.Range("K132").Value = 38495.39999999999
$ws.Range("L132").Value = 106032.15
$ws.Range("M132").Value = -35965.39999999999
$ws.Range("N132").Value = -111092.15
$ws.Range("H134").Value = 69928.5
$ws.Range("J134").Value = 69928.5
$ws.Range("L134").Value = 209785.5
$ws.Range("N134").Value = -214855.5
